# Skill information translation (~170)
# Applies Korean translations to the skill table (rows 141-170, columns A/B/D)
# and fixes a mistranslation in D130 ("눈사람 무효" -> "오뚝이 무효").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D130").Value = "북풍의 사냥꾼(北風の狩人) + 풍압【대】 무효(風圧【大】無効) + 오뚝이 무효(だるま無効)"
$ws.Range("A141").Value = "차지단축(溜め短縮)"
$ws.Range("B141").Value = "집중(集中)/잡념(雑念)"
$ws.Range("D141").Value = "대검, 해머, 활의 차지 시간이 0.8배가 되고 태도, 쌍검, 슬래시 액스, 챠지 액스의 게이지 충전량이 1.2배가 된다./대검, 해머, 활의 차지 시간이 1.2배가 되고 태도, 쌍검, 슬래시 액스, 챠지 액스의 게이지 충전량이 0.8배가 된다."
$ws.Range("A142").Value = "오뚝이(だるま)"
$ws.Range("B142").Value = "오뚝이 무효(だるま無効)"
$ws.Range("D142").Value = "눈사람 상태, 거품 상태【대】, 뼈 투성이 상태를 무효화"
$ws.Range("A143").Value = "탄약절약(弾薬節約)"
$ws.Range("B143").Value = "탄약절약(弾薬節約)"
$ws.Range("D143").Value = "사격 시에 20%의 확률로 탄이나 병을 소모하지 않는다."
$ws.Range("A144").Value = "찬스(チャンス)"
$ws.Range("B144").Value = "비장의 카드(切り札)"
$ws.Range("D144").Value = "같은 에어리어의 대형 몬스터가 당함, 다운했을 때에 1분간 수기의 게이지와 효과시간이 1.2배가 된다. 그 이외에도 일부 수기의 효과가 상승한다."
$ws.Range("A145").Value = "청각보호(聴覚保護)"
$ws.Range("B145").Value = "고급 귀마개(高級耳栓)/귀마개(耳栓)"
$ws.Range("D145").Value = "포효【소】와 포효【대】를 무효화. 데미지가 있는 포효는 막을 수 없다./포효【소】를 무효화"
$ws.Range("A146").Value = "조합수(調合数)"
$ws.Range("B146").Value = "최대수생산(最大数生産)"
$ws.Range("D146").Value = "조합의 생산수가 반드시 최대가 된다."
$ws.Range("A147").Value = "조합성공률(調合成功率)"
$ws.Range("B147").Value = "조합성공률+45%(調合成功率+45%)/조합성공률+20%(調合成功率+20%)/조합성공률-10%(調合成功率-10%)/조합성공률-20%(調合成功率-20%)"
$ws.Range("D147").Value = "조합 성공률+45%/조합 성공률+20%/조합 성공률-10%/조합 성공률-20%"
$ws.Range("A148").Value = "도약(跳躍)"
$ws.Range("B148").Value = "비연(飛燕)"
$ws.Range("D148").Value = "점프 공격시에 위력과 단차 속성축적치, 다운 속성축적치가 1.1배로 상승"
$ws.Range("A149").Value = "통격(痛撃)"
$ws.Range("B149").Value = "약점특효(弱点特効)"
$ws.Range("D149").Value = "공격시에 육질이 45%이상이면 회심률에 50%의 보정이 붙는다."
$ws.Range("A150").Value = "통상탄 강화(通常弾強化)"
$ws.Range("B150").Value = "통상탄・연사 화살UP(通常弾・連射矢UP)"
$ws.Range("D150").Value = "통상탄, 연사 화살의 위력이 1.1배가 된다."
$ws.Range("A151").Value = "통상탄추가(通常弾追加)"
$ws.Range("B151").Value = "통상탄 전LV 추가(通常弾全LV追加)"
$ws.Range("D151").Value = "전LV의 통상탄을 사용할 수 있게 된다."
$ws.Range("A152").Value = "천안(天眼)"
$ws.Range("B152").Value = "천안의 혼(天眼の魂)"
$ws.Range("D152").Value = "간파+3(見切り+3) + 도전자+1(挑戦者+1)"
$ws.Range("A153").Value = "몸통 배가(胴系統倍加)"
$ws.Range("B153").Value = "몸통 배가(胴系統倍加)"
$ws.Range("D153").Value = "몸 파츠의 스킬 포인트가 2배가 된다."
$ws.Range("A154").Value = "투혼(闘魂)"
$ws.Range("B154").Value = "도전자+2(挑戦者+2)/도전자+1(挑戦者+1)"
$ws.Range("D154").Value = "대형 몬스터의 분노에 반응해서 3초후에 공격력이 20, 회심률이 15% 상승한다. 스킬 힘의해방(力の解放), 풀 챠지(フルチャージ)와는 중복되지 않는다./대형 몬스터의 분노에 반응해서 3초후에 공격력이 10, 회심률이 10% 상승한다. 스킬 힘의해방(力の解放), 풀 챠지(フルチャージ)와는 중복되지 않는다."
$ws.Range("A155").Value = "도장(刀匠)"
$ws.Range("B155").Value = "진타(真打)"
$ws.Range("D155").Value = "예리도 레벨+1(斬れ味レベル+1) + 공격력UP【대】(攻撃力UP【大】)"
$ws.Range("A156").Value = "연마사(研ぎ師)"
$ws.Range("B156").Value = "숫돌 사용 고속화(砥石使用高速化)/숫돌 사용 저속화(砥石使用低速化)"
$ws.Range("D156").Value = "숫돌 계열 아이템의 연마 시간이 짧아진다. 통상의 4회에서 1회가 된다./숫돌 계열 아이템의 연마 시간이 길어진다. 통상의 4회에서 5회가 된다."
$ws.Range("A157").Value = "독(毒)"
$ws.Range("B157").Value = "독내성(毒耐性)/독배가(毒倍加)"
$ws.Range("D157").Value = "독, 맹독상태가 되지 않는다. 극독은 맹독으로 경감된다./독에 의해 받는 데미지가 2배가 된다."
$ws.Range("A158").Value = "특수회심(特殊会心)"
$ws.Range("B158").Value = "회심격【특수】(会心撃【特殊】)"
$ws.Range("D158").Value = "크리티컬 공격 시에 가하는 상태이상치(마비, 독, 수면)가 1.2배가 된다."
$ws.Range("A159").Value = "특수공격(特殊攻撃)"
$ws.Range("B159").Value = "상태이상공격+2(状態異常攻撃+2)/상태이상공격+1(状態異常攻撃+1)/상태이상공격약화(状態異常攻撃弱化)"
$ws.Range("D159").Value = "독, 마비, 수면, 포획마취의 속성치가 1.2배+1이 된다./독, 마비, 수면, 포획마취의 속성치가 1.1배+1이 된다./독, 마비, 수면, 포획마취의 속성치가 0.9배가 된다."
$ws.Range("A160").Value = "독병추가(毒瓶追加)"
$ws.Range("B160").Value = "독병 추가(毒ビン追加)"
$ws.Range("D160").Value = "독병의 장착이 가능하게 된다."
$ws.Range("A161").Value = "둔기(鈍器)"
$ws.Range("B161").Value = "둔기사용(鈍器使い)"
$ws.Range("D161").Value = "예리도가 나쁠수록 공격력 상승. 게이지가 녹색：1.1배, 황색：1.15배, 귤색이나 적색：1.2배"
$ws.Range("A162").Value = "육식(肉食)"
$ws.Range("B162").Value = "고기 애호가(お肉大好き)"
$ws.Range("D162").Value = "생고기를 먹을 수 있게 되며 스태미너의 최대치가 50 상승한다. 또, 덜 익힌 고기나 탄 고기, 쿨 미트, 핫 미트로 1분간 강주효과가 발생."
$ws.Range("A163").Value = "도난 무효(盗み無効)"
$ws.Range("B163").Value = "도난 무효(盗み無効)"
$ws.Range("D163").Value = "아이템을 도난당하지 않게 된다."
$ws.Range("A164").Value = "납도(納刀)"
$ws.Range("B164").Value = "납도술(納刀術)"
$ws.Range("D164").Value = "납도 속도가 1.4배가 된다."
$ws.Range("A165").Value = "납도연마(納刀研磨)"
$ws.Range("B165").Value = "도전자의 납도(挑戦者の納刀)"
$ws.Range("D165").Value = "몬스터가 분노 상태일 때 납도하면, 25%의 확률로 예리도가 회복. 회복량은 무기종류에 따라 다르다."
$ws.Range("A166").Value = "탑승(乗り)"
$ws.Range("B166").Value = "탑승 명인(乗り名人)/탑승 하수(乗り下手)"
$ws.Range("D166").Value = "탑승 축적치가 1.25배로 늘고, 탑승상태 시 공격 게이지의 증가량이 1.25배가 된다./탑승 축적치가 0.9배로 줄고, 탑승상태 시 공격 게이지의 증가량이 0.9배가 된다."
$ws.Range("A167").Value = "갈무리(剥ぎ取り)"
$ws.Range("B167").Value = "갈무리 명인(剥ぎ取り名人)/갈무리 철인(剥ぎ取り鉄人)"
$ws.Range("D167").Value = "갈무리 회수가 1회 증가하고, 갈무리 중에 움츠리지 않게 된다./갈무리 중에 움츠리지 않게 된다."
$ws.Range("A168").Value = "폭탄강화(爆弾強化)"
$ws.Range("B168").Value = "보머(ボマー)"
$ws.Range("D168").Value = "폭탄의 데미지가 1.3배가 된다. 폭파속성의 축적치가 1.2배가 된다(고양이 화약술과 병용해서 1.25배 상한). 폭탄 조합성공률이 100%가 된다."
$ws.Range("A169").Value = "폭파병추가(爆破瓶追加)"
$ws.Range("B169").Value = "폭파병 추가(爆破ビン追加)"
$ws.Range("D169").Value = "폭파병을 장착 가능하게 된다."
$ws.Range("A170").Value = "벌꿀(ハチミツ)"
$ws.Range("B170").Value = "허니 헌터(ハニーハンター)"
$ws.Range("D170").Value = "벌꿀(ハチミツ)이나 로열 허니(ロイヤルハニー)가 한 번에 2개 채취 가능하다."

# Restore the active cell / selection shown in the sheet view.
$ws.Range("D170").Select()
